$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# repull data, push all data, mean calculation
# Update column F (dSF) values for affected rows
$ws.Range("F8").Value = 0
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = -3
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = -6
$ws.Range("F20").Value = -2
